$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text (daily conversion summary) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 7.25 = 29246.38 pesos"), "✅ 1000 Bs = 7.27 = 29367.27 pesos"
$text = $text -replace [regex]::Escape("✅ 29246.38 pesos = 7.24 = 942.32 Bs"), "✅ 29367.27 pesos = 7.25 = 932.77 Bs"
$cell.Value2 = $text

# --- Update tasas sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 137.5
$wsTasas.Range("O10").Value2 = 4038
$wsTasas.Range("N12").Value2 = 4052
$wsTasas.Range("O12").Value2 = 128.7
